$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.54%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.46%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.132"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.02%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05609"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.44%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.470"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.19%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8195"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.26%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8334"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.89%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1329"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.58%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06941"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.40%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.03090"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.87%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.02893"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.09381"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.04%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.001520"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.65%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'One"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.0005960"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.18%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.006162"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.24%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.653"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'3.38%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'GateToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.025"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.190"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'8.31%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.3112"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.12%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-1.51%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.742"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.12%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04589"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.58%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1342"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.38%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001227"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.52%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004495"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-2.47%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009599"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.99%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03639"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.45%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.1370"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'30.26%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'KickToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.006149"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.12%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002600"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.07%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'6.17%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005356"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.75%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'4.89%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002501"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'17.78%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E50").Style = "Normal"
